# Update 相談件数 sheet with two additional days of data (2020-04-11, 2020-04-12)
# and push the footnote row down from row 78 to row 79.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Row 78 used to hold the footnote text in column B; capture that text
# before the row gets overwritten with real data.
$footnote = $ws.Range("B78").Value2

# Row 77 used to be blank (only A77/D77/E77 carried formatting); fill it
# with real data.
$ws.Range("A77").Value = 43932
$ws.Range("B77").Value = 709
$ws.Range("C77").Value = 20261
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 4872

# Row 78 becomes a real data row too.
$ws.Range("A78").Value = 43933
$ws.Range("B78").Value = 697
$ws.Range("C78").Value = 20958
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 4872

# The formatted-but-empty A/D/E cells that used to live on row 78 move
# down to the new row 79, carrying the same number formatting.
$ws.Range("A77").Copy()
$ws.Range("A79").PasteSpecial(-4122)
$ws.Range("D77").Copy()
$ws.Range("D79").PasteSpecial(-4122)
$ws.Range("E77").Copy()
$ws.Range("E79").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The footnote itself also moves down to row 79.
$ws.Range("B79").Value = $footnote

# Extend the print area to cover the new row.
$wb.Names.Item(1).RefersTo = "=相談件数!`$A`$1:`$E`$83"

# Move the selected cell like the committed change did.
$ws.Range("E78").Select()
